$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (sheet1) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# M26: 0 -> 11532.42 (new June sale of PORCELANATO recorded for this client)
$ws1.Cells.Item(26, 13).Value = 11532.42

# M32: "5 de 30" -> "6 de 30" (count of advisors who sold PORCELANATO)
$ws1.Cells.Item(32, 13).Value = "6 de 30"

# Column F width 13 -> 14 (ColumnWidth excludes the ~0.83 cell padding)
$ws1.Columns.Item(6).ColumnWidth = 13.17

# --- Sheet "VENTA MENSUAL" (sheet2) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# F26 (junio): 0 -> 11532.42
$ws2.Cells.Item(26, 6).Value = 11532.42

# F32 (junio total): 17133.12 -> 28665.54
$ws2.Cells.Item(32, 6).Value = 28665.54

# Column D width 13 -> 14
$ws2.Columns.Item(4).ColumnWidth = 13.17

# --- Sheet "CUMPLIMIENTO MENSUAL" (sheet3) ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 16 (PORCELANATO) VENTA / POR CUMPLIR / CUMPLIMIENTO
$ws3.Cells.Item(16, 4).Value = 19072.49
$ws3.Cells.Item(16, 5).Value = -273.880000000001
$ws3.Cells.Item(16, 6).Value = 1.014569162294446

# Row 19 (TOTAL) VENTA / POR CUMPLIR / CUMPLIMIENTO
$ws3.Cells.Item(19, 4).Value = 28659.78
$ws3.Cells.Item(19, 5).Value = 878.0110755578733
$ws3.Cells.Item(19, 6).Value = 0.9702749920157565
